$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.090.76"
$ws.Range("E2").Value = "  -2.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.896.76"
$ws.Range("E3").Value = "  -2.87%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "332.16"
$ws.Range("E5").Value = "  -3.12%  "
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4596"
$ws.Range("E7").Value = "  -3.75%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4111"
$ws.Range("E8").Value = "  -0.83%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.84"
$ws.Range("E9").Value = "  -0.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07995"
$ws.Range("E10").Value = "  -3.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.010"
$ws.Range("E11").Value = "  -2.52%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.13"
$ws.Range("E12").Value = "  -2.68%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.896.81"
$ws.Range("E13").Value = "  -2.88%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.932"
$ws.Range("E14").Value = "  -4.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.104"
$ws.Range("E15").Value = "  -4.29%  "
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "89.04"
$ws.Range("E17").Value = "  -3.39%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001028"
$ws.Range("E18").Value = "  -3.10%  "
$ws.Range("E19").Value = "  -1.98%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.59"
$ws.Range("E20").Value = "  -2.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "29.053.95"
$ws.Range("E22").Value = "  -2.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.490"
$ws.Range("E23").Value = "  -1.70%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.39"
$ws.Range("E24").Value = "  +1.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.192"
$ws.Range("E25").Value = "  -3.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.120.71"
$ws.Range("E26").Value = "  -2.72%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "157.20"
$ws.Range("E27").Value = "  -2.87%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.74"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.119"
$ws.Range("E29").Value = "  -3.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.601"
$ws.Range("E30").Value = "  -1.78%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "117.01"
$ws.Range("E31").Value = "  -4.80%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.048"
$ws.Range("E32").Value = "  +3.57%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09383"
$ws.Range("E33").Value = "  -2.57%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.417"
$ws.Range("E34").Value = "  -4.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.545"
$ws.Range("E35").Value = "  -3.92%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.349"
$ws.Range("E36").Value = "  -3.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06086"
$ws.Range("E37").Value = "  -3.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02237"
$ws.Range("E38").Value = "  -3.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.420"
$ws.Range("E39").Value = "  -0.92%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.179"
$ws.Range("E40").Value = "  -0.89%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5835"
$ws.Range("E41").Value = "  -4.43%  "
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1826"
$ws.Range("E43").Value = "  -3.60%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.10"
$ws.Range("E44").Value = "  -5.89%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.349"
$ws.Range("E45").Value = "  -1.68%  "
$ws.Range("E46").Value = "  -0.88%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.07499"
$ws.Range("E47").Value = "  +2.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "12.15"
$ws.Range("E48").Value = "  -2.62%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.5540"
$ws.Range("E49").Value = "  -2.93%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.920"
$ws.Range("E50").Value = "  -3.54%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "112.26"
$ws.Range("E51").Value = "  -1.18%  "
